# PUBLONS.xlsx update: expand "Test Cases" sheet with new PUBLONS034-041 rows,
# clear old row 32/33 detail cells, relocate their data to new rows 40/41,
# and append 5 blank rows at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# ---------------------------------------------------------------------------
# Capture the data currently sitting in B32:D33 before we touch anything,
# since that content moves down to rows 40/41.
# ---------------------------------------------------------------------------
$b32 = $ws.Range("B32").Value
$c32 = $ws.Range("C32").Value
$d32 = $ws.Range("D32").Value
$b33 = $ws.Range("B33").Value
$c33 = $ws.Range("C33").Value
$d33 = $ws.Range("D33").Value

# ---------------------------------------------------------------------------
# Establish a reusable "bordered + wrap text" blank-cell style (matches the
# new cellXfs entry s=12) by borrowing it from another sheet that already
# has a plain bordered cell, then turning on wrap text.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$ws.Range("C32").WrapText = $true

# Now C32 carries our target "blank, bordered, wrap" style. Use it as the
# copy source for every blank C-column cell we create below.
$ws.Range("C32").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C34:C39").PasteSpecial(-4122)
$ws.Range("C40:C41").PasteSpecial(-4122)
$ws.Range("C42:C46").PasteSpecial(-4122)

# Clear the values (PasteSpecial(Formats) only copied formatting, cells are
# already blank, but make sure no stray content remains).
$ws.Range("C32:C46").ClearContents()
$ws.Range("C32").WrapText = $true
$ws.Range("C33:C46").WrapText = $true

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 32 and 33: blank out B/C/D (border-only style for B/D, bordered+wrap
# for C). A32/A33 keep their existing TCID text/style.
# ---------------------------------------------------------------------------
$ws.Range("B32").Value = ""
$ws.Range("B32").Borders.LineStyle = 1
$ws.Range("D32").Value = ""
$ws.Range("D32").Borders.LineStyle = 1

$ws.Range("B33").Value = ""
$ws.Range("B33").Borders.LineStyle = 1
$ws.Range("D33").Value = ""
$ws.Range("D33").Borders.LineStyle = 1

Write-Host "done"
